$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 125, pushing the existing rows 125-131
# (the weekly Jengibre / Vega Modelo de Temuco records) down to 126-132.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with this week's record.
$ws.Cells.Item(125, 1).Value = 10
$ws.Cells.Item(125, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(125, 3).Value = "La Araucanía"
$ws.Cells.Item(125, 4).Value = 44610
$ws.Cells.Item(125, 5).Value = 9
$ws.Cells.Item(125, 6).Value = 100114007
$ws.Cells.Item(125, 7).Value = "Jengibre"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 40
$ws.Cells.Item(125, 11).Value = 25000
$ws.Cells.Item(125, 12).Value = 26000
$ws.Cells.Item(125, 13).Value = 25500
$ws.Cells.Item(125, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(125, 15).Value = "Perú"
$ws.Cells.Item(125, 16).Value = 1962
$ws.Cells.Item(125, 17).Value = 13
$ws.Cells.Item(125, 18).Value = "Hortaliza"
